{"js": "// The \"Trigger:\" paragraph holds the sentence\n// \"By clicking on the LogOut option.<nbsp>\" split across several\n// adjacent runs (artifacts of earlier edits). The commit squashes\n// those adjacent runs back into single runs without changing the\n// visible text, for the two groups:\n//   \"By \" + \"clicking on\" + \" the \"      -> \"By clicking on the \"\n//   \" \" + \"option\" + \".<nbsp>\"           -> \" option.<nbsp>\"\n// (the trailing character is a non-breaking space, U+00A0, not a\n// plain space - it must be preserved exactly).\n// Re-typing the exact same text over each found range via\n// insertText(..., Replace) merges the run boundaries while leaving\n// the surrounding runs (incl. the LogOut proofErr-wrapped run) intact.\n\nconst body = context.document.body;\n\nconst firstGroupText = \"By clicking on the \";\nconst firstGroup = body.search(firstGroupText, { matchCase: true });\nfirstGroup.load(\"items\");\nawait context.sync();\n\nif (firstGroup.items.length > 0) {\n  firstGroup.items[0].insertText(firstGroupText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst secondGroupText = \" option.\\u00a0\";\nconst secondGroup = body.search(secondGroupText, { matchCase: true });\nsecondGroup.load(\"items\");\nawait context.sync();\n\nif (secondGroup.items.length > 0) {\n  secondGroup.items[0].insertText(secondGroupText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The \"Trigger:\" paragraph holds the sentence\n#   \"By clicking on the LogOut option.<nbsp>\"\n# split across several adjacent runs (artifacts of earlier edits). The\n# commit squashes those adjacent runs back into single runs without\n# changing the visible text, for the two groups:\n#   \"By \" + \"clicking on\" + \" the \"      -> \"By clicking on the \"\n#   \" \" + \"option\" + \".<nbsp>\"           -> \" option.<nbsp>\"\n# (the trailing character is a non-breaking space, U+00A0, not a plain\n# space - it is preserved because we only replace the matched text).\n#\n# Using Find/Replace on the *middle* run of each group (\"clicking on\",\n# \"option\") and replacing it with its own text merges it together with\n# its immediate same-formatting neighbours on both sides into a single\n# run, which is exactly the collapsing the diff shows - while leaving\n# the unrelated, differently-formatted run of leading spaces, and the\n# LogOut run (wrapped in proofErr spell-check tags), untouched.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.Execute(\"clicking on\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"clicking on\", $wdReplaceAll) | Out-Null\n\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Execute(\"option\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"option\", $wdReplaceAll) | Out-Null\n"}
